# Update column G (K) values for rows 2-12 on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = 3
    4  = 1
    5  = 2
    6  = 2
    7  = 3
    8  = 0
    9  = 4
    10 = 4
    11 = 2
    12 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
